# Add a new row (40) to the "Translation" sheet for the new "Autotest" screen text entry.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$ws.Range("B40").Value = "SingleUseId43"
$ws.Range("C40").Value = "Typography_00"
$ws.Range("D40").Value = "Left"
$ws.Range("E40").Value = "LTR"
$ws.Range("F40").Value = "Autotest:"
